# Update shared-string labels used across the workbook's sheets
# "Selfreport without App" -> "Selfreport w/o App"
# "Selfreport with App"    -> "Selfreport w/ App"
# "Sensor + Selfreport without App" -> "Sensor + Selfreport w/o App"
# "Sensor + Selfreport with App"    -> "Sensor + Selfreport w/ App"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        switch ($cell.Value2) {
            "Selfreport without App" { $cell.Value2 = "Selfreport w/o App" }
            "Selfreport with App" { $cell.Value2 = "Selfreport w/ App" }
            "Sensor + Selfreport without App" { $cell.Value2 = "Sensor + Selfreport w/o App" }
            "Sensor + Selfreport with App" { $cell.Value2 = "Sensor + Selfreport w/ App" }
        }
    }
}
